$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's price report was added to the dataset. It belongs right
# after the existing "Primera" / "Segunda" pair for 2021-01-19 (rows
# 122-123), so insert a fresh row at 124 and push everything else down
# by one (old row 124 -> 125, ..., old row 142 -> 143).
$ws.Rows("124").Insert()

# Populate the newly inserted row 124 with the new weekly record.
$ws.Range("A124").Value = 11
$ws.Range("B124").Value = "Vega Monumental Concepción"
$ws.Range("C124").Value = "Bíobío"
$ws.Range("D124").Value = 44504
$ws.Range("E124").Value = 8
$ws.Range("F124").Value = 100114013
$ws.Range("G124").Value = "Zanahoria"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 350
$ws.Range("K124").Value = 6000
$ws.Range("L124").Value = 7000
$ws.Range("M124").Value = 6571
$ws.Range("N124").Value = "$/saco 20 kilos"
$ws.Range("O124").Value = "Chillán"
$ws.Range("P124").Value = 329
$ws.Range("Q124").Value = 20
$ws.Range("R124").Value = "Hortaliza"

# Keep column D's date format on the new row consistent with the rest
# of the column.
$ws.Range("D124").NumberFormat = "YYYY-MM-DD HH:MM:SS"
